# Scheduled market-data refresh: updates the computed pricing/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ -> columns H-N) on each "leve" sheet with freshly pulled
# Universalis market data. Columns A-G (leve name/item/level/exp/gil/amount/
# item id) are static reference data and are left untouched.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1651.5
$ws.Range("I4").Value = 815.625
$ws.Range("K4").Value = 815.625
$ws.Range("M4").Value = -701.625
$ws.Range("H11").Value = 90860.53999999999
$ws.Range("I11").Value = 90860.53999999999
$ws.Range("K11").Value = 90860.53999999999
$ws.Range("M11").Value = -90720.53999999999
$ws.Range("H17").Value = 259532.23
$ws.Range("J17").Value = 265309.4
$ws.Range("L17").Value = 795928.2000000001
$ws.Range("N17").Value = -796264.2000000001
$ws.Range("H53").Value = 84456.836
$ws.Range("I53").Value = 236.16667
$ws.Range("J53").Value = 168677.5
$ws.Range("K53").Value = 236.16667
$ws.Range("L53").Value = 168677.5
$ws.Range("M53").Value = 400.83333
$ws.Range("N53").Value = -169951.5
$ws.Range("H70").Value = 1462882.8
$ws.Range("J70").Value = 6016.3335
$ws.Range("L70").Value = 18049.0005
$ws.Range("N70").Value = -18589.0005
$ws.Range("H73").Value = 1462882.8
$ws.Range("J73").Value = 6016.3335
$ws.Range("L73").Value = 18049.0005
$ws.Range("N73").Value = -19921.0005
$ws.Range("H86").Value = 25132400
$ws.Range("I86").Value = 4798.75
$ws.Range("K86").Value = 4798.75
$ws.Range("M86").Value = -3675.75
$ws.Range("H88").Value = 1138
$ws.Range("I88").Value = 985.4286
$ws.Range("K88").Value = 985.4286
$ws.Range("M88").Value = -579.4286
$ws.Range("H89").Value = 25132400
$ws.Range("I89").Value = 4798.75
$ws.Range("K89").Value = 23993.75
$ws.Range("M89").Value = -18377.75
$ws.Range("H91").Value = 1138
$ws.Range("I91").Value = 985.4286
$ws.Range("K91").Value = 985.4286
$ws.Range("M91").Value = 418.5714
$ws.Range("H98").Value = 781.8570999999999
$ws.Range("I98").Value = 679.9474
$ws.Range("J98").Value = 1750
$ws.Range("K98").Value = 679.9474
$ws.Range("L98").Value = 1750
$ws.Range("M98").Value = 818.0526
$ws.Range("N98").Value = -4746
$ws.Range("H112").Value = 49471.668
$ws.Range("I112").Value = 336083.34
$ws.Range("J112").Value = 1703.0555
$ws.Range("K112").Value = 1008250.02
$ws.Range("L112").Value = 5109.166499999999
$ws.Range("M112").Value = -1007142.02
$ws.Range("N112").Value = -7325.166499999999
$ws.Range("H116").Value = 64928930
$ws.Range("I116").Value = 41835664
$ws.Range("K116").Value = 41835664
$ws.Range("M116").Value = -41832222
$ws.Range("H122").Value = 781.8570999999999
$ws.Range("I122").Value = 679.9474
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 2039.8422
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = 410.1578
$ws.Range("N122").Value = -10150
$ws.Range("H129").Value = 1246.7142
$ws.Range("I129").Value = 537.7143
$ws.Range("K129").Value = 1613.1429
$ws.Range("M129").Value = 3386.8571
$ws.Range("H132").Value = 2198.7544
$ws.Range("I132").Value = 2256.9424
$ws.Range("K132").Value = 6770.8272
$ws.Range("M132").Value = -4240.8272
$ws.Range("H138").Value = 1297.5294
$ws.Range("I138").Value = 1071.8334
$ws.Range("J138").Value = 1839.2
$ws.Range("K138").Value = 3215.5002
$ws.Range("L138").Value = 5517.6
$ws.Range("M138").Value = 1924.4998
$ws.Range("N138").Value = -15797.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 620097.7
$ws.Range("I45").Value = 1390494.8
$ws.Range("K45").Value = 1390494.8
$ws.Range("M45").Value = -1390117.8
$ws.Range("H74").Value = 1360.4147
$ws.Range("I74").Value = 1312.8718
$ws.Range("K74").Value = 1312.8718
$ws.Range("M74").Value = -438.8717999999999
$ws.Range("H77").Value = 1360.4147
$ws.Range("I77").Value = 1312.8718
$ws.Range("K77").Value = 6564.358999999999
$ws.Range("M77").Value = -2196.358999999999
$ws.Range("H88").Value = 23812980
$ws.Range("I88").Value = 41668340
$ws.Range("K88").Value = 41668340
$ws.Range("M88").Value = -41667934
$ws.Range("H91").Value = 23812980
$ws.Range("I91").Value = 41668340
$ws.Range("K91").Value = 41668340
$ws.Range("M91").Value = -41666936

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 12492.75
$ws.Range("I26").Value = 12492.75
$ws.Range("K26").Value = 12492.75
$ws.Range("M26").Value = -12200.75
$ws.Range("H96").Value = 14464.929
$ws.Range("I96").Value = 14464.929
$ws.Range("K96").Value = 14464.929
$ws.Range("M96").Value = -11718.929
$ws.Range("H107").Value = 45456228
$ws.Range("I107").Value = 1615
$ws.Range("K107").Value = 1615
$ws.Range("M107").Value = 305

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 976.6429000000001
$ws.Range("J22").Value = 899.8
$ws.Range("L22").Value = 899.8
$ws.Range("N22").Value = -1599.8
$ws.Range("H31").Value = 2813
$ws.Range("I31").Value = 1227.1305
$ws.Range("J31").Value = 3572.8958
$ws.Range("K31").Value = 1227.1305
$ws.Range("L31").Value = 3572.8958
$ws.Range("M31").Value = -932.1305
$ws.Range("N31").Value = -4162.8958
$ws.Range("H34").Value = 2813
$ws.Range("I34").Value = 1227.1305
$ws.Range("J34").Value = 3572.8958
$ws.Range("K34").Value = 1227.1305
$ws.Range("L34").Value = 3572.8958
$ws.Range("M34").Value = -1025.1305
$ws.Range("N34").Value = -3976.8958
$ws.Range("H62").Value = 6949.6
$ws.Range("J62").Value = 6666
$ws.Range("L62").Value = 6666
$ws.Range("N62").Value = -7914
$ws.Range("H65").Value = 6949.6
$ws.Range("J65").Value = 6666
$ws.Range("L65").Value = 33330
$ws.Range("N65").Value = -39570
$ws.Range("H99").Value = 507071.2
$ws.Range("I99").Value = 618976
$ws.Range("J99").Value = 3499.5
$ws.Range("K99").Value = 618976
$ws.Range("L99").Value = 3499.5
$ws.Range("M99").Value = -617478
$ws.Range("N99").Value = -6495.5
$ws.Range("H107").Value = 2904.5454
$ws.Range("I107").Value = 2767.125
$ws.Range("K107").Value = 2767.125
$ws.Range("M107").Value = -847.125
$ws.Range("H126").Value = 507071.2
$ws.Range("I126").Value = 618976
$ws.Range("J126").Value = 3499.5
$ws.Range("K126").Value = 1856928
$ws.Range("L126").Value = 10498.5
$ws.Range("M126").Value = -1854458
$ws.Range("N126").Value = -15438.5
$ws.Range("H134").Value = 747.58826
$ws.Range("I134").Value = 721.0714
$ws.Range("K134").Value = 2163.2142
$ws.Range("M134").Value = 371.7857999999997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 32.5625
$ws.Range("I2").Value = 34.066666
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 34.066666
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 78.933334
# N2 did not previously exist on this row (no HQ listings priced); the refreshed
# market pull now yields an HQ profit figure, so the cell is newly populated.
$ws.Range("N2").Value = -236
$ws.Range("H102").Value = 1408.4667
$ws.Range("I102").Value = 1276.9642
$ws.Range("J102").Value = 3249.5
$ws.Range("K102").Value = 1276.9642
$ws.Range("L102").Value = 3249.5
$ws.Range("M102").Value = 345.0358000000001
$ws.Range("N102").Value = -6493.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26318790
$ws.Range("I7").Value = 38463616
$ws.Range("J7").Value = 5000.8335
$ws.Range("K7").Value = 38463616
$ws.Range("L7").Value = 5000.8335
$ws.Range("M7").Value = -38463504
$ws.Range("N7").Value = -5224.8335
$ws.Range("H22").Value = 734
$ws.Range("I22").Value = 391
$ws.Range("J22").Value = 962.6667
$ws.Range("K22").Value = 391
$ws.Range("L22").Value = 962.6667
$ws.Range("M22").Value = -96
$ws.Range("N22").Value = -1552.6667
$ws.Range("H27").Value = 734
$ws.Range("I27").Value = 391
$ws.Range("J27").Value = 962.6667
$ws.Range("K27").Value = 391
$ws.Range("L27").Value = 962.6667
$ws.Range("M27").Value = -284
$ws.Range("N27").Value = -1176.6667
$ws.Range("H40").Value = 2924
$ws.Range("I40").Value = 2471.3333
$ws.Range("J40").Value = 3700
$ws.Range("K40").Value = 2471.3333
$ws.Range("L40").Value = 3700
$ws.Range("M40").Value = -2335.3333
$ws.Range("N40").Value = -3972
$ws.Range("H55").Value = 454.76923
$ws.Range("J55").Value = 1048.5
$ws.Range("L55").Value = 1048.5
$ws.Range("N55").Value = -1394.5
$ws.Range("H126").Value = 26318790
$ws.Range("I126").Value = 38463616
$ws.Range("J126").Value = 5000.8335
$ws.Range("K126").Value = 115390848
$ws.Range("L126").Value = 15002.5005
$ws.Range("M126").Value = -115388378
$ws.Range("N126").Value = -19942.5005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 5000
$ws.Range("K58").Value = 5000
$ws.Range("M58").Value = -4692
$ws.Range("H126").Value = 1399
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
# HQ average price/listing price collapsed to 0 on refresh (no HQ market data),
# so the derived HQ-profit figure (N126) is no longer meaningful and is cleared.
$ws.Range("N126").ClearContents()
